$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 204. This shifts the existing rows 204-283 down to 205-284,
# matching the target dimension change from A1:R283 to A1:R284.
$ws.Rows("204").Insert()

# Populate the newly inserted row 204. Columns A, B, C, E, F, G, H, I, N, O, Q, R
# keep the same values as the row that used to be row 204 (now row 205), while
# D, J, K, L, M, P get new values per the edit.
$ws.Range("A204").Value = 10
$ws.Range("B204").Value = "Vega Modelo de Temuco"
$ws.Range("C204").Value = "La Araucanía"
$ws.Range("D204").Value = 44636
$ws.Range("E204").Value = 9
$ws.Range("F204").Value = 100114013
$ws.Range("G204").Value = "Zanahoria"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 30
$ws.Range("K204").Value = 7000
$ws.Range("L204").Value = 7000
$ws.Range("M204").Value = 7000
$ws.Range("N204").Value = "$/saco 25 kilos"
$ws.Range("O204").Value = "Región de La Araucanía"
$ws.Range("P204").Value = 280
$ws.Range("Q204").Value = 25
$ws.Range("R204").Value = "Hortaliza"
